# "moved previous concerts list" - the workbook file was relocated into the
# repo's _data/ folder and re-opened; the author then widened the data
# columns (C:G) to fit their content and left the selection back at the
# top of the sheet (A2) instead of the old scroll position (F60 / A23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns C:G to fit the concert data (jour, mois, annee, libelle,
# programme, distribution, lieu, manif, pdf_affiche columns live in A:I;
# C:G are the ones that received explicit custom widths).
$ws.Columns("C").ColumnWidth = 12.666666666666666
$ws.Columns("D").ColumnWidth = 23.333333333333336
$ws.Columns("E").ColumnWidth = 103.5
$ws.Columns("F").ColumnWidth = 74.83333333333334
$ws.Columns("G").ColumnWidth = 49.83333333333333

# Reset the view back to the top-left of the sheet with A2 selected
# (previously it was scrolled down to A23 with F60 selected).
$ws.Range("A2").Select()
